$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ReimbUnits")

# Rename the "Create Custom flow1" label to "Create ReimbUnit1"
$ws.Range("D2").Value = "Create ReimbUnit1"

# Remove the second test case (AutoUnit2) row entirely, shifting
# everything below it up by one row
$ws.Rows("3").Delete()

$ws.Range("A8").Select()
